$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '255.17'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.72'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '6.133'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06028'

$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.445'
$ws.Range("E6").Value = '5GateTokenGT'

$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.608'
$ws.Range("E7").Value = '6KuCoinTokenKCS'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.323'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8031'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1529'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07984'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03354'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03113'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09301'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.599'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001691'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04795'

$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0006123'
$ws.Range("E18").Value = '17OneONEWorstin24h'

$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006291'
$ws.Range("E19").Value = '18TigerCashTCH'

$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.005894'
$ws.Range("E20").Value = '19HotbitTokenHTB'

$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001095'
$ws.Range("E21").Value = '20BitKanKAN'

$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001509'
$ws.Range("E22").Value = '21NitroExNTX'

$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.665'
$ws.Range("E23").Value = '22LEOLEO'

$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.204'
$ws.Range("E24").Value = '23BTSETokenBTSE'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3349'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0006520'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04488'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007071'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1073'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003371'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01075'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.002477'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00005916'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000755'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.7049'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.09412'
$ws.Range("E49").Value = '48BOLOBOLO'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002115'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.01017'
